$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort paths by sorting indices then reshuffle instead of sorting in place
# -> add a new benchmark row ("110* - sort indices and shuffle") below the
#    existing table, which pushes the chart further down/right.

$ws.Range("A10").Value = "110* - sort indices and shuffle"
$ws.Range("B10").Value = "ON"
$ws.Range("C10").Value = "ON"
$ws.Range("D10").Value = "OFF"
$ws.Range("E10").Value = 510.15899999999999
$ws.Range("F10").Value = 9.8008699999999997

# Re-home the chart so it sits below the newly added row, matching the
# new two-cell anchor recorded for the drawing.
$co = $ws.ChartObjects().Item(1)
$co.Left = 201.1875
$co.Top = 234.37496062992125
$co.Width = 309.0938287401575
$co.Height = 311.62503937007875

# Update the active selection / view to reflect where the user ended up.
$ws.Range("C14").Select()
